$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Create the new "features" worksheet right after Sheet1
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "features"

$basic = @(
"1. Introduction to Frappe framework and the ERPNext application",
"2. Installation of Frappe Framework",
"3. Creatingasite to work with ERPNext",
"4. Installing ERPNext onasite",
"5. Frappe Field Types",
"6. Introduction to the frappe desk",
"7. Creatingacustom app",
"8. Installing the custom app on your site",
"9. Exploring the Frappe DocTypes",
"10. Working with Dependent Fields",
"11. Introduction to Custom Scripts",
"12. Using Custom Scripts to make API Calls",
"13. Using Custom Scripts to do Arithmetics",
"14. Working with Server Scripts",
"15. Working with Menus in Custom Applications",
"16. Email Functionality",
"17. Creating custom reports",
"18. Customizing the print format",
"19. Working with schedulers",
"20. Introduction to the Frappe Dialog API",
"21. URL Routing",
"22. Working with Currencies",
"23. Implementing content pagination",
"24. Exploring the Frappe Chart"
)

$advanced = @(
"1. Working with Controllers, ORM and SQL in Frappe Framework",
"2. Hooks and DocType manipulation from events",
"3. ERPNext integration with devices (We will integrateabiometric device)",
"4. Writing custom APIS within the ERPNext ecosystem",
"5. Consuming custom APIS from other software applications",
"6. Sending data from ERPNext to other applications with Web Hooks",
"7. Using Events to manipulate Doctypes",
"8. Building web apps in Frappe Framework",
"9. Working withabootstrap template on Frappe Framework",
"10. Working with dynamic data on the custom web app",
"11. Error Handling in Frappe Framework",
"12. Building custom pages in Frappe Framework",
"13. Managing your ERPNext instance and custom apps with GitHub",
"14. Setting up Frappe for multi-tenancy",
"15. ERPNext for E-Commerce"
)

# Rows 2-25: Basic feature list (entered first so these strings take the
# earlier shared-string slots, matching the authored workbook)
$row = 2
foreach ($item in $basic) {
    $ws.Cells.Item($row, 1).Value = $item
    $row++
}

# Row 1: "Basic" header (bold)
$ws.Range("A1").Value = "Basic"
$ws.Range("A1").Font.Bold = $true

# Row 26 stays blank, row 27: "Advanced" header (bold)
$ws.Range("A27").Value = "Advanced"
$ws.Range("A27").Font.Bold = $true

# Rows 28-42: Advanced feature list
$row = 28
foreach ($item in $advanced) {
    $ws.Cells.Item($row, 1).Value = $item
    $row++
}

# Match column A width to the authored bestFit width (~67.4 characters)
$ws.Columns("A:A").ColumnWidth = 66.6

# Portrait page orientation (mirrors Sheet1's print setup)
$ws.PageSetup.Orientation = 1

# Leave the features sheet with A36 selected
$ws.Range("A36").Select()

# Re-activate Sheet1 and restore its selection
$sheet1.Activate()
$sheet1.Range("B29").Select()
